$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for price cells that would otherwise be parsed as numbers
$textCells = @("D5", "D6", "D7", "D10", "D13", "D20", "D21", "D24", "D26", "D29", "D32", "D33", "D39", "D42", "D43", "D45", "D46", "D47", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Update Price (D) and Volume(1h) (E) columns with latest values
$ws.Range("D2").Value = "55.811.66"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").Value = "2.375.83"
$ws.Range("E3").Value = "  -4.55%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "478.74"
$ws.Range("E5").Value = "  -1.81%  "
$ws.Range("D6").Value = "147.22"
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("E8").Value = "  -2.21%  "
$ws.Range("D9").Value = "2.376.55"
$ws.Range("E9").Value = "  -5.30%  "
$ws.Range("D10").Value = "0.0972"
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("E11").Value = "  -6.36%  "
$ws.Range("E12").Value = "  -2.75%  "
$ws.Range("D13").Value = "0.125"
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("D14").Value = "2.787.35"
$ws.Range("E14").Value = "  -4.80%  "
$ws.Range("D15").Value = "55.890.35"
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("E16").Value = "  -4.33%  "
$ws.Range("E17").Value = "  -3.30%  "
$ws.Range("D18").Value = "2.369.58"
$ws.Range("E18").Value = "  -5.61%  "
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").Value = "315.44"
$ws.Range("E20").Value = "  -1.59%  "
$ws.Range("D21").Value = "9.69"
$ws.Range("E21").Value = "  -4.98%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("E23").Value = "  -2.95%  "
$ws.Range("D24").Value = "56.75"
$ws.Range("E24").Value = "  -3.17%  "
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("D26").Value = "0.396"
$ws.Range("E26").Value = "  -3.93%  "
$ws.Range("E27").Value = "  -5.96%  "
$ws.Range("D28").Value = "2.483.19"
$ws.Range("E28").Value = "  -4.94%  "
$ws.Range("D29").Value = "7.25"
$ws.Range("E29").Value = "  -5.03%  "
$ws.Range("D30").Value = "0.0₃0771"
$ws.Range("E30").Value = "  -2.44%  "
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").Value = "147.28"
$ws.Range("E32").Value = "  -1.21%  "
$ws.Range("D33").Value = "18.02"
$ws.Range("E33").Value = "  -1.36%  "
$ws.Range("E34").Value = "  -1.48%  "
$ws.Range("E35").Value = "  -2.91%  "
$ws.Range("E36").Value = "  -3.70%  "
$ws.Range("E37").Value = "  -4.17%  "
$ws.Range("E38").Value = "  -3.48%  "
$ws.Range("D39").Value = "33.42"
$ws.Range("E39").Value = "  -2.26%  "
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("D42").Value = "3.38"
$ws.Range("E42").Value = "  -4.24%  "
$ws.Range("D43").Value = "0.0537"
$ws.Range("E43").Value = "  -3.82%  "
$ws.Range("E44").Value = "  +3.60%  "
$ws.Range("D45").Value = "0.582"
$ws.Range("E45").Value = "  -5.71%  "
$ws.Range("D46").Value = "10.19"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").Value = "254.76"
$ws.Range("E47").Value = "  -2.23%  "
$ws.Range("E48").Value = "  -2.29%  "
$ws.Range("D49").Value = "4.52"
$ws.Range("E49").Value = "  -7.12%  "
$ws.Range("D50").Value = "16.87"
$ws.Range("E50").Value = "  -4.53%  "
$ws.Range("D51").Value = "1.778.10"
$ws.Range("E51").Value = "  -7.02%  "
